$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for new columns I (9) and J (10) ---
$ws.Columns.Item(9).ColumnWidth = 18.43
$ws.Columns.Item(10).ColumnWidth = 44.3

# --- Header row ---
$ws.Range("J1").Value = "Links"
$ws.Range("I1").Value = "Papers"
$ws.Range("I1").WrapText = $true
$ws.Range("J1").WrapText = $true

# --- Row 2 ---
$ws.Range("I2").Value = "Chuck & Norris (2014);Du und Ich (2012)"
$ws.Range("J2").Value = "https://example.com;https://example.com"
$ws.Hyperlinks.Add($ws.Range("J2"), "https://example.com")

# --- Row 3 (wrap) ---
$ws.Range("I3").Value = "Atkinson & Piketty (2010);Gottfried & Schellhorn (2004)"
$ws.Range("I3").WrapText = $true
$ws.Range("J3").Value = "https://ideas.repec.org/p/iaw/iawdip/15.html;https://ideas.repec.org/p/iaw/iawdip/15.html`n"
$ws.Range("J3").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("J3"), "https://ideas.repec.org/p/iaw/iawdip/15.html")
$ws.Rows.Item(3).RowHeight = 105

# --- Row 4 (wrap) ---
$ws.Range("I4").Value = "Doerrenberg et. al. (2017);Clementi and Gallegati (2005)"
$ws.Range("J4").Value = "https://ideas.repec.org/a/eee/pubeco/v151y2017icp41-55.html;https://ideas.repec.org/p/arx/papers/physics-0504217.html`n"
$ws.Range("J4").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("J4"), "https://ideas.repec.org/a/eee/pubeco/v151y2017icp41-55.html")
$ws.Rows.Item(4).RowHeight = 75

# --- Row 5 ---
$ws.Range("I5").Value = "Peter et. al. (2018)"
$ws.Range("J5").Value = "https://ideas.repec.org/p/diw/diwwpp/dp1770.html"
$ws.Hyperlinks.Add($ws.Range("J5"), "https://ideas.repec.org/p/diw/diwwpp/dp1770.html")

# --- Row 6 ---
$ws.Range("J6").Value = "https://ideas.repec.org/p/ces/ceswps/_8382.html"
$ws.Hyperlinks.Add($ws.Range("J6"), "https://ideas.repec.org/p/ces/ceswps/_8382.html")
$ws.Range("I6").Value = "Falk et. al. (2020)"
$ws.Range("I6").WrapText = $true

# --- Row 7 ---
$ws.Range("J7").Value = "https://ideas.repec.org/a/uwp/jhriss/v54y2019i2p468-502.html"
$ws.Range("I7").Value = "Marcus & Zambre (2019)"

# --- Row 8 ---
$ws.Range("J8").Value = "https://ideas.repec.org/p/arx/papers/1909.08299.html;https://ideas.repec.org/a/eee/ecoedu/v41y2014icp14-23.html"
$ws.Hyperlinks.Add($ws.Range("J8"), "https://ideas.repec.org/p/arx/papers/1909.08299.html")
$ws.Hyperlinks.Add($ws.Range("J7"), "https://ideas.repec.org/a/uwp/jhriss/v54y2019i2p468-502.html")
$ws.Range("I8").Value = "Gorgen & Schienle (2019);Bruckmeier & Wigger (2014)"

# --- Rows 9-12 (shared text, only row 9 gets a hyperlink ref but no hyperlink registered) ---
$ws.Range("I9").Value = "Lechner et. al. (2011)"
$ws.Range("J9").Value = "https://ideas.repec.org/a/bla/jeurec/v9y2011i4p742-784.html"
$ws.Range("I10").Value = "Lechner et. al. (2011)"
$ws.Range("J10").Value = "https://ideas.repec.org/a/bla/jeurec/v9y2011i4p742-784.html"
$ws.Range("I11").Value = "Lechner et. al. (2011)"
$ws.Range("J11").Value = "https://ideas.repec.org/a/bla/jeurec/v9y2011i4p742-784.html"
$ws.Range("I12").Value = "Lechner et. al. (2011)"
$ws.Range("J12").Value = "https://ideas.repec.org/a/bla/jeurec/v9y2011i4p742-784.html"

# --- Rows 13-14 ---
$ws.Range("I13").Value = "Biewen et. al. (2014)"
$ws.Range("J13").Value = "https://ideas.repec.org/a/ucp/jlabec/doi10.1086-677233.html"
$ws.Hyperlinks.Add($ws.Range("J13"), "https://ideas.repec.org/a/ucp/jlabec/doi10.1086-677233.html")
$ws.Range("I14").Value = "Biewen et. al. (2014)"
$ws.Range("J14").Value = "https://ideas.repec.org/a/ucp/jlabec/doi10.1086-677233.html"
$ws.Hyperlinks.Add($ws.Range("J14"), "https://ideas.repec.org/a/ucp/jlabec/doi10.1086-677233.html")

# --- Rows 15-16 (hyperlink style, but NO hyperlink relationship) ---
$ws.Range("I15").Value = "Caliendo & Künn (2011) "
$ws.Range("J15").Value = "https://ideas.repec.org/a/eee/pubeco/v95y2011i3p311-331.html"
$ws.Range("J15").Style = "Hyperlink"
$ws.Range("I16").Value = "Caliendo & Künn (2011) "
$ws.Range("J16").Value = "https://ideas.repec.org/a/eee/pubeco/v95y2011i3p311-331.html"
$ws.Range("J16").Style = "Hyperlink"

# --- Row 17 ---
$ws.Range("I17").Value = "Caliendo et. al. (2016)"
$ws.Range("J17").Value = "https://ideas.repec.org/a/eee/eecrev/v86y2016icp87-108.html"
$ws.Hyperlinks.Add($ws.Range("J17"), "https://ideas.repec.org/a/eee/eecrev/v86y2016icp87-108.html")

# --- Row 18 ---
$ws.Range("I18").Value = "Doerr et. al. (2014);Huber et al. (2018)"
$ws.Range("J18").Value = "https://ideas.repec.org/p/iza/izadps/dp8454.html;https://ideas.repec.org/p/cpr/ceprdp/10650.html"
$ws.Hyperlinks.Add($ws.Range("J18"), "https://ideas.repec.org/p/iza/izadps/dp8454.html")

# --- Rows 19-20 ---
$ws.Range("J19").Value = "https://ideas.repec.org/p/zbw/fubsbe/201817.html"
$ws.Hyperlinks.Add($ws.Range("J19"), "https://ideas.repec.org/p/zbw/fubsbe/201817.html")
$ws.Range("I19").Value = "Thiedig (2018)"
$ws.Range("I20").Value = "Thiedig (2018)"
$ws.Range("J20").Value = "https://ideas.repec.org/p/zbw/fubsbe/201817.html"
$ws.Hyperlinks.Add($ws.Range("J20"), "https://ideas.repec.org/p/zbw/fubsbe/201817.html")

# --- Rows 21-23 ---
$ws.Range("I21").Value = "Hohmeyer & Wolff (2010)"
$ws.Range("J21").Value = "https://ideas.repec.org/p/iab/iabdpa/201021.html"
$ws.Hyperlinks.Add($ws.Range("J21"), "https://ideas.repec.org/p/iab/iabdpa/201021.html")
$ws.Range("I22").Value = "Hohmeyer & Wolff (2010)"
$ws.Range("J22").Value = "https://ideas.repec.org/p/iab/iabdpa/201021.html"
$ws.Hyperlinks.Add($ws.Range("J22"), "https://ideas.repec.org/p/iab/iabdpa/201021.html")
$ws.Range("I23").Value = "Hohmeyer & Wolff (2010)"
$ws.Range("J23").Value = "https://ideas.repec.org/p/iab/iabdpa/201021.html"
$ws.Hyperlinks.Add($ws.Range("J23"), "https://ideas.repec.org/p/iab/iabdpa/201021.html")

# --- Selection: set active cell to J23 ---
$ws.Range("J23").Select()
